$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Operators")
$ws2 = $wb.Worksheets.Item("Tasks")

# --- Operators sheet: columns G (Not evening) and H (Not task) ---
$ws1.Range("G1").Value = "Not evening"
$ws1.Range("H1").Value = "Not task"
$ws1.Range("H2").Value = "26,27"
$ws1.Range("H3").Value = "8,12,13,26,27"
$ws1.Range("G4").Value = "2,9,18,23,30"
$ws1.Range("H4").Value = "5,6,19,20"
$ws1.Range("G5").Value = "4,9,16"
$ws1.Range("H5").Value = "5,6,26,27"
$ws1.Range("G6").Value = "9,25"
$ws1.Range("H6").Value = "26,27"
$ws1.Range("G7").Value = "9,18,23,30"
$ws1.Range("H7").Value = "5,6,19,20"
$ws1.Range("G8").Value = "2,9,15,23"
$ws1.Range("H8").Value = "26,27"
$ws1.Range("G9").Value = "2,9"
$ws1.Range("H9").Value = "5,6,12,13"
$ws1.Range("G10").Value = "4,9,15,22,29"
$ws1.Range("H10").Value = "5,6,26,27"
$ws1.Range("G11").Value = "2,11,17,23"
$ws1.Range("H11").Value = "1,3,12,13,24"
$ws1.Range("G12").Value = 23
$ws1.Range("H12").Value = "12,13,26,27"
$ws1.Range("G13").Value = "11,14,23"
$ws1.Range("H13").Value = "5,6,26,27"
$ws1.Range("G14").Value = 23
$ws1.Range("H14").Value = "26,27,28"
$ws1.Range("H15").Value = "11,12,13"
$ws1.Range("G16").Value = "14,25,31"
$ws1.Range("H16").Value = "9,12,13,26,27"
$ws1.Range("H17").Value = "2,3,4,5,6,8,10,11,17,23,25,26,27"
$ws1.Range("G18").Value = 8
$ws1.Range("H18").Value = "9,10,11,14,15,16,17,18,19,20,21,23,25,26,27,28"
$ws1.Range("G19").Value = "2,9"
$ws1.Range("H19").Value = "11,12,13,14,26,27"
$ws1.Range("G20").Value = "3,23,29"
$ws1.Range("H20").Value = "7,17,26,27"
$ws1.Range("G22").Value = 10
$ws1.Range("H22").Value = "9,17,26,27"
$ws1.Range("G23").Value = "10,18,23"
$ws1.Range("H23").Value = "5,6,19,20"
$ws1.Range("H24").Value = "2,3,5,6,10,21,26,27"
$ws1.Range("G25").Value = "17,25"
$ws1.Range("H25").Value = "3,10,21,23,26,27,28"
$ws1.Range("H26").Value = "8,10,11,12,13,25,26,27"
$ws1.Range("H27").Value = "19,20,26,27,30"
$ws1.Range("G28").Value = 10
$ws1.Range("H28").Value = "19,20,26,27"
$ws1.Range("G29").Value = 7
$ws1.Range("H29").Value = "3,5,6,17,25,26,27"
$ws1.Range("G30").Value = "7,15,23"
$ws1.Range("H30").Value = "10,11,26,27,28"
$ws1.Range("G31").Value = 7
$ws1.Range("H31").Value = "10,11,18,19,20,21,26,27"
$ws1.Range("G32").Value = "7,23"
$ws1.Range("H32").Value = "10,11,26,27,28"
$ws1.Range("H33").Value = "1,26,27"
$ws1.Range("G34").Value = 7
$ws1.Range("H34").Value = "3,14,17"
$ws1.Range("H35").Value = "1,3,4,5,6,9,16,23,26,27"
$ws1.Range("G38").Value = 25
$ws1.Range("H38").Value = "9,19,20,26,27,28"
$ws1.Range("H40").Value = "26,27,28"
$ws1.Range("G41").Value = 18
$ws1.Range("H42").Value = "1,11"
$ws1.Range("H43").Value = "2,9,16,23"

# --- Tasks sheet: column G (probability) ---
$ws2.Range("G1").Value = "probability"
$ws2.Range("G2").Value = 1
$ws2.Range("G3").Value = 1
$ws2.Range("G4").Value = 1
$ws2.Range("G5").Value = 1
$ws2.Range("G6").Value = 0
$ws2.Range("G7").Value = 0
$ws2.Range("G8").Value = 1
$ws2.Range("G9").Value = 1
$ws2.Range("G10").Value = 1
$ws2.Range("G11").Value = 1
$ws2.Range("G12").Value = 1
$ws2.Range("G13").Value = 0.25
$ws2.Range("G14").Value = 0.2
$ws2.Range("G15").Value = 0
$ws2.Range("G16").Value = 0
$ws2.Range("G17").Value = 1
$ws2.Range("G18").Value = 0
$ws2.Range("G19").Value = 0.4

$ws1.Range("D11").Select() | Out-Null
$ws2.Range("H12").Select() | Out-Null
$ws2.Activate() | Out-Null
